$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.31647966666667
$ws.Range("H2").Value = 54.949439
$ws.Range("I2").Value = 0.005487334033884006
$ws.Range("J2").Value = 0.005487334033884005
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 298.1092115270585
$ws.Range("R2").Value = 2682.982903743527
$ws.Range("S2").Value = 0.0003637312578819696
$ws.Range("T2").Value = 0.0003637312578819695
$ws.Range("G3").Value = 18.31647966666667
$ws.Range("H3").Value = 54.949439
$ws.Range("I3").Value = 0.005487334033884006
$ws.Range("J3").Value = 0.005487334033884005
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 1565.384568191207
$ws.Range("R3").Value = 14088.46111372087
$ws.Range("S3").Value = 0.001909968816933155
$ws.Range("T3").Value = 0.001909968816933155
$ws.Range("G4").Value = 18.31647966666667
$ws.Range("H4").Value = 54.949439
$ws.Range("I4").Value = 0.005487334033884006
$ws.Range("J4").Value = 0.005487334033884005
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 780.5295138622695
$ws.Range("R4").Value = 7024.765624760425
$ws.Range("S4").Value = 0.0009523455529496658
$ws.Range("T4").Value = 0.0009523455529496658
$ws.Range("G5").Value = 18.31647966666667
$ws.Range("H5").Value = 54.949439
$ws.Range("I5").Value = 0.005487334033884006
$ws.Range("J5").Value = 0.005487334033884005
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 1853.321344194802
$ws.Range("R5").Value = 16679.89209775322
$ws.Range("S5").Value = 0.002261288406119215
$ws.Range("T5").Value = 0.002261288406119215
$ws.Range("I6").Value = 0.9472399998689139
$ws.Range("J6").Value = 0.9472399998689137
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 51460.50299546646
$ws.Range("R6").Value = 463144.5269591981
$ws.Range("S6").Value = 0.06278837674924016
$ws.Range("T6").Value = 0.06278837674924015
$ws.Range("I7").Value = 0.9472399998689139
$ws.Range("J7").Value = 0.9472399998689137
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.3297045251354632
$ws.Range("T7").Value = 0.3297045251354631
$ws.Range("I8").Value = 0.9472399998689139
$ws.Range("J8").Value = 0.9472399998689137
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 134737.337301345
$ws.Range("R8").Value = 1212636.035712105
$ws.Range("S8").Value = 0.1643967354421623
$ws.Range("T8").Value = 0.1643967354421623
$ws.Range("I9").Value = 0.9472399998689139
$ws.Range("J9").Value = 0.9472399998689137
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 319925.8690999618
$ws.Range("R9").Value = 2879332.821899656
$ws.Range("S9").Value = 0.3903503625420483
$ws.Range("T9").Value = 0.3903503625420482
$ws.Range("G10").Value = 155.6514383333333
$ws.Range("H10").Value = 466.954315
$ws.Range("I10").Value = 0.04663076369111781
$ws.Range("J10").Value = 0.0466307636911178
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 2533.299433026199
$ws.Range("R10").Value = 22799.69489723579
$ws.Range("S10").Value = 0.003090948396549844
$ws.Range("T10").Value = 0.003090948396549843
$ws.Range("G11").Value = 155.6514383333333
$ws.Range("H11").Value = 466.954315
$ws.Range("I11").Value = 0.04663076369111781
$ws.Range("J11").Value = 0.0466307636911178
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 13302.46663212151
$ws.Range("R11").Value = 119722.1996890935
$ws.Range("S11").Value = 0.01623070584182637
$ws.Range("T11").Value = 0.01623070584182637
$ws.Range("G12").Value = 155.6514383333333
$ws.Range("H12").Value = 466.954315
$ws.Range("I12").Value = 0.04663076369111781
$ws.Range("J12").Value = 0.0466307636911178
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 6632.854331467134
$ws.Range("R12").Value = 59695.68898320421
$ws.Range("S12").Value = 0.00809292821571677
$ws.Range("T12").Value = 0.008092928215716768
$ws.Range("G13").Value = 155.6514383333333
$ws.Range("H13").Value = 466.954315
$ws.Range("I13").Value = 0.04663076369111781
$ws.Range("J13").Value = 0.0466307636911178
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 15749.32182207289
$ws.Range("R13").Value = 141743.896398656
$ws.Range("S13").Value = 0.01921618123702483
$ws.Range("T13").Value = 0.01921618123702482
$ws.Range("G14").Value = 2.142642
$ws.Range("H14").Value = 6.427926
$ws.Range("I14").Value = 0.0006419024060843985
$ws.Range("J14").Value = 0.0006419024060843984
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 34.872493450102
$ws.Range("R14").Value = 313.852441050918
$ws.Range("S14").Value = 0.00004254888952646482
$ws.Range("T14").Value = 0.00004254888952646481
$ws.Range("G15").Value = 2.142642
$ws.Range("H15").Value = 6.427926
$ws.Range("I15").Value = 0.0006419024060843985
$ws.Range("J15").Value = 0.0006419024060843984
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 183.11699535049
$ws.Range("R15").Value = 1648.05295815441
$ws.Range("S15").Value = 0.0002234260884365693
$ws.Range("T15").Value = 0.0002234260884365693
$ws.Range("G16").Value = 2.142642
$ws.Range("H16").Value = 6.427926
$ws.Range("I16").Value = 0.0006419024060843985
$ws.Range("J16").Value = 0.0006419024060843984
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 91.30549915027599
$ws.Range("R16").Value = 821.749492352484
$ws.Range("S16").Value = 0.0001114043537512646
$ws.Range("T16").Value = 0.0001114043537512646
$ws.Range("G17").Value = 2.142642
$ws.Range("H17").Value = 6.427926
$ws.Range("I17").Value = 0.0006419024060843985
$ws.Range("J17").Value = 0.0006419024060843984
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 216.799528284624
$ws.Range("R17").Value = 1951.195754561616
$ws.Range("S17").Value = 0.0002645230743700998
$ws.Range("T17").Value = 0.0002645230743700998
